$d = $word.ActiveDocument

function New-PkgXml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1) Add a new centered "Kartikey Mishra" byline paragraph right after the
#    "Pawat Saengsiripongpun" paragraph (and before the blank spacer lines).
# ---------------------------------------------------------------------------
$pawat = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Pawat Saengsiripongpun") {
        $pawat = $p
        break
    }
}

$pawat.Range.InsertParagraphAfter() | Out-Null
$kartikeyPara = $pawat.Next()
$kartikeyXml = New-PkgXml '<w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Kartikey</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Mishra</w:t></w:r></w:p></w:body>'
$kartikeyPara.Range.InsertXML($kartikeyXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge "Spacebar to " + "jump" into one run (drop the gramStart/gramEnd
#    proofErr wrapper) and append a new run about the jump sound.
# ---------------------------------------------------------------------------
$spacebarPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Spacebar to jump") {
        $spacebarPara = $p
        break
    }
}

$spacebarXml = New-PkgXml '<w:body><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Spacebar to jump</w:t></w:r><w:r><w:t xml:space="preserve"> (contain sound when jumping)</w:t></w:r></w:p></w:body>'
$spacebarPara.Range.InsertXML($spacebarXml) | Out-Null

# ---------------------------------------------------------------------------
# 3) Re-emit "How to play" heading with its formatting on the paragraph mark
#    too, then insert a new paragraph describing the background music.
# ---------------------------------------------------------------------------
$howToPlayPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "How to play") {
        $howToPlayPara = $p
        break
    }
}

$howToPlayXml = New-PkgXml '<w:body><w:p><w:pPr><w:rPr><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="23"/><w:szCs w:val="23"/><w:u w:val="single"/></w:rPr><w:t>How to play</w:t></w:r></w:p></w:body>'
$howToPlayPara.Range.InsertXML($howToPlayXml) | Out-Null

$howToPlayPara.Range.InsertParagraphAfter() | Out-Null
$bgMusicPara = $howToPlayPara.Next()
$bgMusicXml = New-PkgXml '<w:body><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/></w:rPr><w:t>Background music automatically starts at the beginning of the game</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>.</w:t></w:r></w:p></w:body>'
$bgMusicPara.Range.InsertXML($bgMusicXml) | Out-Null
